$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (header "K") values recalculated after switching from Strike# to K
$newValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 7
    6  = 2
    7  = 2
    8  = 4
    9  = 2
    10 = 1
    11 = 6
    12 = 2
    13 = 6
    14 = 5
    15 = 1
    16 = 1
    17 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
